# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates metrics for row 3 (file_name = metrics_sim_with_priors.json)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.6
$ws.Range("D3").Value = 0.9333333333333333
$ws.Range("F3").Value = 0.9333333333333333
$ws.Range("H3").Value = 0.6627906976744186
$ws.Range("I3").Value = 0.11440329218107
$ws.Range("J3").Value = 0.5333333333333333
$ws.Range("K3").Value = 35.8

$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 5
$ws.Range("S3").Value = 13
$ws.Range("T3").Value = 21
$ws.Range("U3").Value = 36
$ws.Range("V3").Value = 240
$ws.Range("W3").Value = 238
$ws.Range("X3").Value = 230
$ws.Range("Y3").Value = 222
$ws.Range("Z3").Value = 207

$ws.Range("AF3").Value = 0.987654
$ws.Range("AG3").Value = 0.979424
$ws.Range("AH3").Value = 0.946502
$ws.Range("AI3").Value = 0.9135799999999999
$ws.Range("AJ3").Value = 0.8518520000000001
